$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year column (O) mirrors the formatting of the existing N column
$ws.Range("N3:N5").Copy()
$ws.Range("O3:O5").PasteSpecial(-4122)

$ws.Cells.Item(3, 15).Value = 2021
$ws.Cells.Item(4, 15).Value = 14
$ws.Cells.Item(5, 15).Value = 1252.8

$ws.Range("O9").Select()
